$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "personal_account_iban" column (D) with header + value,
# matching the existing header/value style used by the other columns.
$ws.Range("D1").Value = "personal_account_iban"
$ws.Range("D2").Value = "205-9031004417882-84"

# Give the new column a sensible custom width (mirrors the A/B columns
# which already carry an explicit width).
$ws.Columns.Item(4).ColumnWidth = 27.6675

# Move the active selection to D4, matching where the author's cursor
# ended up after the edit.
$ws.Range("D4").Select() | Out-Null
